$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the publication Date value
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Update the Contact value
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new "Jurisdiction" row right after "Contact", pushing the rest down
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
